$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 180, shifting existing rows 180-288 down to 181-289.
$ws.Rows("180:180").Insert()

# Populate the newly inserted row 180 with the new record's data.
$ws.Cells.Item(180, 1).Value2 = 3
$ws.Cells.Item(180, 2).Value2 = "Femacal de La Calera"
$ws.Cells.Item(180, 3).Value2 = "Coquimbo"
$ws.Cells.Item(180, 4).Value2 = 44606
$ws.Cells.Item(180, 5).Value2 = 5
$ws.Cells.Item(180, 6).Value2 = 100112009
$ws.Cells.Item(180, 7).Value2 = "Acelga"
$ws.Cells.Item(180, 8).Value2 = "Sin especificar"
$ws.Cells.Item(180, 9).Value2 = "Primera"
$ws.Cells.Item(180, 10).Value2 = 230
$ws.Cells.Item(180, 11).Value2 = 2500
$ws.Cells.Item(180, 12).Value2 = 2800
$ws.Cells.Item(180, 13).Value2 = 2657
$ws.Cells.Item(180, 14).Value2 = "$/docena de atados (6 kilos)"
$ws.Cells.Item(180, 15).Value2 = "Provincia de Quillota"
$ws.Cells.Item(180, 16).Value2 = 443
$ws.Cells.Item(180, 17).Value2 = 6
$ws.Cells.Item(180, 18).Value2 = "Hortaliza"
